$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 1.624282666666667
$ws.Range("H2").Value = 4.872847999999999
$ws.Range("I2").Value = 0.08561172663893989
$ws.Range("J2").Value = 0.08561172663893987
$ws.Range("K2").Value = 1.0
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.04936366666666667
$ws.Range("N2").Value = 0.148091
$ws.Range("O2").Value = 0.04616170608573571
$ws.Range("P2").Value = 0.0461617060857357
$ws.Range("Q2").Value = 0.08018054812977778
$ws.Range("R2").Value = 0.7216249331679999
$ws.Range("S2").Value = 0.003951983362599093
$ws.Range("T2").Value = 0.003951983362599091
# Row 3
$ws.Range("G3").Value = 1.624282666666667
$ws.Range("H3").Value = 4.872847999999999
$ws.Range("I3").Value = 0.08561172663893989
$ws.Range("J3").Value = 0.08561172663893987
$ws.Range("O3").Value = 0.008057748967298944
$ws.Range("P3").Value = 0.008057748967298944
$ws.Range("Q3").Value = 0.01399590231111111
$ws.Range("R3").Value = 0.1259631208
$ws.Range("S3").Value = 0.0006898378019135974
$ws.Range("T3").Value = 0.0006898378019135973
# Row 4
$ws.Range("G4").Value = 1.624282666666667
$ws.Range("H4").Value = 4.872847999999999
$ws.Range("I4").Value = 0.08561172663893989
$ws.Range("J4").Value = 0.08561172663893987
$ws.Range("M4").Value = 1.011383666666666
$ws.Range("N4").Value = 3.034151
$ws.Range("O4").Value = 0.9457805449469654
$ws.Range("P4").Value = 0.9457805449469653
$ws.Range("Q4").Value = 1.642772959116444
$ws.Range("R4").Value = 14.784956632048
$ws.Range("S4").Value = 0.0809699054744272
$ws.Range("T4").Value = 0.08096990547442717
# Row 5
$ws.Range("I5").Value = 0.1920894545885022
$ws.Range("J5").Value = 0.1920894545885022
$ws.Range("K5").Value = 1.0
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.04936366666666667
$ws.Range("N5").Value = 0.148091
$ws.Range("O5").Value = 0.04616170608573571
$ws.Range("P5").Value = 0.0461617060857357
$ws.Range("Q5").Value = 0.1799033656196667
$ws.Range("R5").Value = 1.619130290577
$ws.Range("S5").Value = 0.008867176944883712
$ws.Range("T5").Value = 0.008867176944883712
# Row 6
$ws.Range("I6").Value = 0.1920894545885022
$ws.Range("J6").Value = 0.1920894545885022
$ws.Range("O6").Value = 0.008057748967298944
$ws.Range("P6").Value = 0.008057748967298944
$ws.Range("S6").Value = 0.001547808604339521
$ws.Range("T6").Value = 0.001547808604339521
# Row 7
$ws.Range("I7").Value = 0.1920894545885022
$ws.Range("J7").Value = 0.1920894545885022
$ws.Range("M7").Value = 1.011383666666666
$ws.Range("N7").Value = 3.034151
$ws.Range("O7").Value = 0.9457805449469654
$ws.Range("P7").Value = 0.9457805449469653
$ws.Range("Q7").Value = 3.685936192599666
$ws.Range("R7").Value = 33.173425733397
$ws.Range("S7").Value = 0.1816744690392789
$ws.Range("T7").Value = 0.1816744690392789
# Row 8
$ws.Range("E8").Value = 3.0
$ws.Range("F8").Value = 1.0
$ws.Range("G8").Value = 2.741573333333333
$ws.Range("H8").Value = 8.22472
$ws.Range("I8").Value = 0.1445012198865677
$ws.Range("J8").Value = 0.1445012198865677
$ws.Range("K8").Value = 1.0
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.04936366666666667
$ws.Range("N8").Value = 0.148091
$ws.Range("O8").Value = 0.04616170608573571
$ws.Range("P8").Value = 0.0461617060857357
$ws.Range("Q8").Value = 0.1353341121688889
$ws.Range("R8").Value = 1.21800700952
$ws.Range("S8").Value = 0.006670422841434005
$ws.Range("T8").Value = 0.006670422841434005
# Row 9
$ws.Range("E9").Value = 3.0
$ws.Range("F9").Value = 1.0
$ws.Range("G9").Value = 2.741573333333333
$ws.Range("H9").Value = 8.22472
$ws.Range("I9").Value = 0.1445012198865677
$ws.Range("J9").Value = 0.1445012198865677
$ws.Range("O9").Value = 0.008057748967298944
$ws.Range("P9").Value = 0.008057748967298944
$ws.Range("Q9").Value = 0.02362322355555556
$ws.Range("R9").Value = 0.212609012
$ws.Range("S9").Value = 0.001164354555314428
$ws.Range("T9").Value = 0.001164354555314428
# Row 10
$ws.Range("E10").Value = 3.0
$ws.Range("F10").Value = 1.0
$ws.Range("G10").Value = 2.741573333333333
$ws.Range("H10").Value = 8.22472
$ws.Range("I10").Value = 0.1445012198865677
$ws.Range("J10").Value = 0.1445012198865677
$ws.Range("M10").Value = 1.011383666666666
$ws.Range("N10").Value = 3.034151
$ws.Range("O10").Value = 0.9457805449469654
$ws.Range("P10").Value = 0.9457805449469653
$ws.Range("Q10").Value = 2.772782490302221
$ws.Range("R10").Value = 24.95504241271999
$ws.Range("S10").Value = 0.1366664424898193
$ws.Range("T10").Value = 0.1366664424898192
# Row 11
$ws.Range("G11").Value = 2.34427
$ws.Range("H11").Value = 7.03281
$ws.Range("I11").Value = 0.1235603916279767
$ws.Range("J11").Value = 0.1235603916279767
$ws.Range("K11").Value = 1.0
$ws.Range("L11").Value = 0.3333333333333333
$ws.Range("M11").Value = 0.04936366666666667
$ws.Range("N11").Value = 0.148091
$ws.Range("O11").Value = 0.04616170608573571
$ws.Range("P11").Value = 0.0461617060857357
$ws.Range("Q11").Value = 0.1157217628566667
$ws.Range("R11").Value = 1.04149586571
$ws.Range("S11").Value = 0.005703758482169058
$ws.Range("T11").Value = 0.005703758482169057
# Row 12
$ws.Range("G12").Value = 2.34427
$ws.Range("H12").Value = 7.03281
$ws.Range("I12").Value = 0.1235603916279767
$ws.Range("J12").Value = 0.1235603916279767
$ws.Range("O12").Value = 0.008057748967298944
$ws.Range("P12").Value = 0.008057748967298944
$ws.Range("Q12").Value = 0.02019979316666666
$ws.Range("R12").Value = 0.1817981385
$ws.Range("S12").Value = 0.000995618618039382
$ws.Range("T12").Value = 0.000995618618039382
# Row 13
$ws.Range("G13").Value = 2.34427
$ws.Range("H13").Value = 7.03281
$ws.Range("I13").Value = 0.1235603916279767
$ws.Range("J13").Value = 0.1235603916279767
$ws.Range("M13").Value = 1.011383666666666
$ws.Range("N13").Value = 3.034151
$ws.Range("O13").Value = 0.9457805449469654
$ws.Range("P13").Value = 0.9457805449469653
$ws.Range("Q13").Value = 2.370956388256666
$ws.Range("R13").Value = 21.33860749431
$ws.Range("S13").Value = 0.1168610145277682
$ws.Range("T13").Value = 0.1168610145277682
# Row 14
$ws.Range("G14").Value = 4.282534999999999
$ws.Range("H14").Value = 12.847605
$ws.Range("I14").Value = 0.225721312715906
$ws.Range("J14").Value = 0.225721312715906
$ws.Range("K14").Value = 1.0
$ws.Range("L14").Value = 0.3333333333333333
$ws.Range("M14").Value = 0.04936366666666667
$ws.Range("N14").Value = 0.148091
$ws.Range("O14").Value = 0.04616170608573571
$ws.Range("P14").Value = 0.0461617060857357
$ws.Range("Q14").Value = 0.2114016302283333
$ws.Range("R14").Value = 1.902614672055
$ws.Range("S14").Value = 0.01041968089487809
$ws.Range("T14").Value = 0.01041968089487809
# Row 15
$ws.Range("G15").Value = 4.282534999999999
$ws.Range("H15").Value = 12.847605
$ws.Range("I15").Value = 0.225721312715906
$ws.Range("J15").Value = 0.225721312715906
$ws.Range("O15").Value = 0.008057748967298944
$ws.Range("P15").Value = 0.008057748967298944
$ws.Range("Q15").Value = 0.03690117658333333
$ws.Range("R15").Value = 0.33211058925
$ws.Range("S15").Value = 0.001818805674433954
$ws.Range("T15").Value = 0.001818805674433954
# Row 16
$ws.Range("G16").Value = 4.282534999999999
$ws.Range("H16").Value = 12.847605
$ws.Range("I16").Value = 0.225721312715906
$ws.Range("J16").Value = 0.225721312715906
$ws.Range("M16").Value = 1.011383666666666
$ws.Range("N16").Value = 3.034151
$ws.Range("O16").Value = 0.9457805449469654
$ws.Range("P16").Value = 0.9457805449469653
$ws.Range("Q16").Value = 4.331285950928332
$ws.Range("R16").Value = 38.98157355835499
$ws.Range("S16").Value = 0.213482826146594
$ws.Range("T16").Value = 0.2134828261465939
# Row 17
$ws.Range("E17").Value = 3.0
$ws.Range("F17").Value = 1.0
$ws.Range("G17").Value = 4.335555666666667
$ws.Range("H17").Value = 13.006667
$ws.Range("I17").Value = 0.2285158945421077
$ws.Range("J17").Value = 0.2285158945421077
$ws.Range("K17").Value = 1.0
$ws.Range("L17").Value = 1.0
$ws.Range("M17").Value = 0.04936366666666667
$ws.Range("N17").Value = 0.148091
$ws.Range("O17").Value = 0.04616170608573571
$ws.Range("P17").Value = 0.0461617060857357
$ws.Range("Q17").Value = 0.2140189247441111
$ws.Range("R17").Value = 1.926170322697
$ws.Range("S17").Value = 0.01054868355977175
$ws.Range("T17").Value = 0.01054868355977175
# Row 18
$ws.Range("E18").Value = 3.0
$ws.Range("F18").Value = 1.0
$ws.Range("G18").Value = 4.335555666666667
$ws.Range("H18").Value = 13.006667
$ws.Range("I18").Value = 0.2285158945421077
$ws.Range("J18").Value = 0.2285158945421077
$ws.Range("O18").Value = 0.008057748967298944
$ws.Range("P18").Value = 0.008057748967298944
$ws.Range("Q18").Value = 0.03735803799444445
$ws.Range("R18").Value = 0.33622234195
$ws.Range("S18").Value = 0.001841323713258063
$ws.Range("T18").Value = 0.001841323713258063
# Row 19
$ws.Range("E19").Value = 3.0
$ws.Range("F19").Value = 1.0
$ws.Range("G19").Value = 4.335555666666667
$ws.Range("H19").Value = 13.006667
$ws.Range("I19").Value = 0.2285158945421077
$ws.Range("J19").Value = 0.2285158945421077
$ws.Range("M19").Value = 1.011383666666666
$ws.Range("N19").Value = 3.034151
$ws.Range("O19").Value = 0.9457805449469654
$ws.Range("P19").Value = 0.9457805449469653
$ws.Range("Q19").Value = 4.384910187190777
$ws.Range("R19").Value = 39.46419168471699
$ws.Range("S19").Value = 0.2161258872690779
$ws.Range("T19").Value = 0.2161258872690779